$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = "x"
$ws.Range("B9").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B17").Select()
